$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2,8).Value = 259.75
$ws.Cells.Item(2,9).Value = 125.42857
$ws.Cells.Item(2,11).Value = 125.42857
$ws.Cells.Item(2,13).Value = -12.42856999999999
$ws.Cells.Item(6,8).Value = 233.70589
$ws.Cells.Item(6,9).Value = 47.727272
$ws.Cells.Item(6,10).Value = 574.6667
$ws.Cells.Item(6,11).Value = 143.181816
$ws.Cells.Item(6,12).Value = 1724.0001
$ws.Cells.Item(6,13).Value = -31.181816
$ws.Cells.Item(6,14).Value = -1948.0001
$ws.Cells.Item(8,8).Value = 310.5
$ws.Cells.Item(8,9).Value = 22
$ws.Cells.Item(8,11).Value = 66
$ws.Cells.Item(8,13).Value = 73
$ws.Cells.Item(9,8).Value = 102.27273
$ws.Cells.Item(9,9).Value = 59.75
$ws.Cells.Item(9,10).Value = 126.57143
$ws.Cells.Item(9,11).Value = 59.75
$ws.Cells.Item(9,12).Value = 126.57143
$ws.Cells.Item(9,13).Value = 109.25
$ws.Cells.Item(9,14).Value = -464.57143
$ws.Cells.Item(11,8).Value = 29.916666
$ws.Cells.Item(11,9).Value = 29.916666
$ws.Cells.Item(11,11).Value = 29.916666
$ws.Cells.Item(11,13).Value = 110.083334
$ws.Cells.Item(38,8).Value = 109
$ws.Cells.Item(38,10).Value = 200
$ws.Cells.Item(38,12).Value = 600
$ws.Cells.Item(38,14).Value = -1344
$ws.Cells.Item(39,8).Value = 215.57895
$ws.Cells.Item(39,9).Value = 49.46154
$ws.Cells.Item(39,10).Value = 575.5
$ws.Cells.Item(39,11).Value = 148.38462
$ws.Cells.Item(39,12).Value = 1726.5
$ws.Cells.Item(39,13).Value = 147.61538
$ws.Cells.Item(39,14).Value = -2318.5
$ws.Cells.Item(43,8).Value = 2471
$ws.Cells.Item(43,10).Value = 2831.6667
$ws.Cells.Item(43,12).Value = 2831.6667
$ws.Cells.Item(43,14).Value = -2969.6667
$ws.Cells.Item(74,8).Value = 6114.5
$ws.Cells.Item(74,9).Value = 5737.4
$ws.Cells.Item(74,11).Value = 5737.4
$ws.Cells.Item(74,13).Value = -4801.4
$ws.Cells.Item(77,8).Value = 6114.5
$ws.Cells.Item(77,9).Value = 5737.4
$ws.Cells.Item(77,11).Value = 28687
$ws.Cells.Item(77,13).Value = -24007
$ws.Cells.Item(112,8).Value = 3305
$ws.Cells.Item(112,10).Value = 3357.5
$ws.Cells.Item(112,12).Value = 10072.5
$ws.Cells.Item(112,14).Value = -12288.5
$ws.Cells.Item(132,8).Value = 12856.348
$ws.Cells.Item(132,9).Value = 12016.105
$ws.Cells.Item(132,10).Value = 16847.5
$ws.Cells.Item(132,11).Value = 36048.315
$ws.Cells.Item(132,12).Value = 50542.5
$ws.Cells.Item(132,13).Value = -33518.315
$ws.Cells.Item(132,14).Value = -55602.5
$ws.Cells.Item(138,8).Value = 2394.8333
$ws.Cells.Item(138,9).Value = 1025.1428
$ws.Cells.Item(138,10).Value = 7188.75
$ws.Cells.Item(138,11).Value = 3075.4284
$ws.Cells.Item(138,12).Value = 21566.25
$ws.Cells.Item(138,13).Value = 2064.5716
$ws.Cells.Item(138,14).Value = -31846.25
$ws.Cells.Item(141,8).Value = 900.6923
$ws.Cells.Item(141,9).Value = 900.6923
$ws.Cells.Item(141,11).Value = 2702.0769
$ws.Cells.Item(141,13).Value = 2477.9231
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32,8).Value = 3032757
$ws.Cells.Item(32,9).Value = 570.8214
$ws.Cells.Item(32,11).Value = 570.8214
$ws.Cells.Item(32,13).Value = -283.8214
$ws.Cells.Item(61,8).Value = 2555.5833
$ws.Cells.Item(61,9).Value = 2376.7
$ws.Cells.Item(61,11).Value = 2376.7
$ws.Cells.Item(61,13).Value = -2164.7
$ws.Cells.Item(74,8).Value = 2878.0417
$ws.Cells.Item(74,10).Value = 4554.5
$ws.Cells.Item(74,12).Value = 4554.5
$ws.Cells.Item(74,14).Value = -6302.5
$ws.Cells.Item(77,8).Value = 2878.0417
$ws.Cells.Item(77,10).Value = 4554.5
$ws.Cells.Item(77,12).Value = 22772.5
$ws.Cells.Item(77,14).Value = -31508.5
$ws.Cells.Item(122,8).Value = 1274.5294
$ws.Cells.Item(122,9).Value = 909.26666
$ws.Cells.Item(122,10).Value = 4014
$ws.Cells.Item(122,11).Value = 2727.79998
$ws.Cells.Item(122,12).Value = 12042
$ws.Cells.Item(122,13).Value = -277.7999799999998
$ws.Cells.Item(122,14).Value = -16942
$ws.Cells.Item(132,8).Value = 2620.75
$ws.Cells.Item(132,9).Value = 2423.7144
$ws.Cells.Item(132,11).Value = 7271.1432
$ws.Cells.Item(132,13).Value = -4741.1432
$ws.Cells.Item(136,8).Value = 2555.5833
$ws.Cells.Item(136,9).Value = 2376.7
$ws.Cells.Item(136,11).Value = 7130.099999999999
$ws.Cells.Item(136,13).Value = -4580.099999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20,8).Value = 1233.3334
$ws.Cells.Item(20,9).Value = 1233.3334
$ws.Cells.Item(20,11).Value = 1233.3334
$ws.Cells.Item(20,13).Value = -986.3334
$ws.Cells.Item(99,8).Value = 1372.1111
$ws.Cells.Item(99,9).Value = 1394.875
$ws.Cells.Item(99,11).Value = 1394.875
$ws.Cells.Item(99,13).Value = 103.125
$ws.Cells.Item(134,8).Value = 6885.077
$ws.Cells.Item(134,9).Value = 1822.8889
$ws.Cells.Item(134,10).Value = 18275
$ws.Cells.Item(134,11).Value = 5468.6667
$ws.Cells.Item(134,12).Value = 54825
$ws.Cells.Item(134,13).Value = -2933.6667
$ws.Cells.Item(134,14).Value = -59895
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22,8).Value = 698.3125
$ws.Cells.Item(22,9).Value = 698.1667
$ws.Cells.Item(22,10).Value = 698.75
$ws.Cells.Item(22,11).Value = 698.1667
$ws.Cells.Item(22,12).Value = 698.75
$ws.Cells.Item(22,13).Value = -348.1667
$ws.Cells.Item(22,14).Value = -1398.75
$ws.Cells.Item(31,8).Value = 4829.7334
$ws.Cells.Item(31,10).Value = 5574.1377
$ws.Cells.Item(31,12).Value = 5574.1377
$ws.Cells.Item(31,14).Value = -6164.1377
$ws.Cells.Item(34,8).Value = 4829.7334
$ws.Cells.Item(34,10).Value = 5574.1377
$ws.Cells.Item(34,12).Value = 5574.1377
$ws.Cells.Item(34,14).Value = -5978.1377
$ws.Cells.Item(134,8).Value = 1812.7037
$ws.Cells.Item(134,9).Value = 906.7619
$ws.Cells.Item(134,10).Value = 4983.5
$ws.Cells.Item(134,11).Value = 2720.2857
$ws.Cells.Item(134,12).Value = 14950.5
$ws.Cells.Item(134,13).Value = -185.2856999999999
$ws.Cells.Item(134,14).Value = -20020.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5,8).Value = 1781.2
$ws.Cells.Item(5,9).Value = 1326.5
$ws.Cells.Item(5,10).Value = 3600
$ws.Cells.Item(5,11).Value = 3979.5
$ws.Cells.Item(5,12).Value = 10800
$ws.Cells.Item(5,13).Value = -3867.5
$ws.Cells.Item(5,14).Value = -11024
$ws.Cells.Item(34,8).Value = 1801.4584
$ws.Cells.Item(34,10).Value = 3635.5454
$ws.Cells.Item(34,12).Value = 10906.6362
$ws.Cells.Item(34,14).Value = -11074.6362
$ws.Cells.Item(55,8).Value = 5330
$ws.Cells.Item(55,10).Value = 5935
$ws.Cells.Item(55,12).Value = 17805
$ws.Cells.Item(55,14).Value = -18159
$ws.Cells.Item(68,8).Value = 574.5
$ws.Cells.Item(68,9).Value = 300
$ws.Cells.Item(68,10).Value = 666
$ws.Cells.Item(68,11).Value = 900
$ws.Cells.Item(68,12).Value = 1998
$ws.Cells.Item(68,13).Value = -89
$ws.Cells.Item(68,14).Value = -3620
$ws.Cells.Item(71,8).Value = 574.5
$ws.Cells.Item(71,9).Value = 300
$ws.Cells.Item(71,10).Value = 666
$ws.Cells.Item(71,11).Value = 2700
$ws.Cells.Item(71,12).Value = 5994
$ws.Cells.Item(71,13).Value = 1356
$ws.Cells.Item(71,14).Value = -14106
$ws.Cells.Item(135,8).Value = 1781.2
$ws.Cells.Item(135,9).Value = 1326.5
$ws.Cells.Item(135,10).Value = 3600
$ws.Cells.Item(135,11).Value = 11938.5
$ws.Cells.Item(135,12).Value = 32400
$ws.Cells.Item(135,13).Value = -9403.5
$ws.Cells.Item(135,14).Value = -37470
$ws.Cells.Item(140,8).Value = 1964.2667
$ws.Cells.Item(140,9).Value = 1964.2667
$ws.Cells.Item(140,10).Value = 0
$ws.Cells.Item(140,11).Value = 5892.800099999999
$ws.Cells.Item(140,12).Value = 0
$ws.Cells.Item(140,13).Value = -712.8000999999995
$ws.Cells.Item(140,14).ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2,8).Value = 112.933334
$ws.Cells.Item(2,10).Value = 92.166664
$ws.Cells.Item(2,12).Value = 92.166664
$ws.Cells.Item(2,14).Value = -318.166664
$ws.Cells.Item(80,8).Value = 2123.75
$ws.Cells.Item(80,9).Value = 1997.5
$ws.Cells.Item(80,11).Value = 1997.5
$ws.Cells.Item(80,13).Value = -999.5
$ws.Cells.Item(83,8).Value = 2123.75
$ws.Cells.Item(83,9).Value = 1997.5
$ws.Cells.Item(83,11).Value = 9987.5
$ws.Cells.Item(83,13).Value = -4995.5
$ws.Cells.Item(122,8).Value = 2681.2222
$ws.Cells.Item(122,9).Value = 1421.4
$ws.Cells.Item(122,11).Value = 4264.200000000001
$ws.Cells.Item(122,13).Value = -1814.200000000001
$ws.Cells.Item(132,8).Value = 21423.754
$ws.Cells.Item(132,9).Value = 24164.457
$ws.Cells.Item(132,10).Value = 3413.4285
$ws.Cells.Item(132,11).Value = 72493.371
$ws.Cells.Item(132,12).Value = 10240.2855
$ws.Cells.Item(132,13).Value = -69963.371
$ws.Cells.Item(132,14).Value = -15300.2855
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40,8).Value = 3763.3333
$ws.Cells.Item(40,9).Value = 3763.3333
$ws.Cells.Item(40,11).Value = 3763.3333
$ws.Cells.Item(40,13).Value = -3627.3333
$ws.Cells.Item(46,8).Value = 997
$ws.Cells.Item(46,10).Value = 1000
$ws.Cells.Item(46,12).Value = 1000
$ws.Cells.Item(46,14).Value = -1376
$ws.Cells.Item(136,8).Value = 4000
$ws.Cells.Item(136,9).Value = 4000
$ws.Cells.Item(136,11).Value = 12000
$ws.Cells.Item(136,13).Value = -9450
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(17,8).Value = 5375
$ws.Cells.Item(17,9).Value = 10000
$ws.Cells.Item(17,10).Value = 750
$ws.Cells.Item(17,11).Value = 10000
$ws.Cells.Item(17,12).Value = 750
$ws.Cells.Item(17,13).Value = -9828
$ws.Cells.Item(17,14).Value = -1094
$ws.Cells.Item(100,8).Value = 1771.4286
$ws.Cells.Item(100,9).Value = 1771.4286
$ws.Cells.Item(100,11).Value = 3542.8572
$ws.Cells.Item(100,13).Value = -3001.8572
$ws.Cells.Item(136,8).Value = 2475.842
$ws.Cells.Item(136,9).Value = 2008.8125
$ws.Cells.Item(136,11).Value = 6026.4375
$ws.Cells.Item(136,13).Value = -3476.4375
